$wb = $excel.ActiveWorkbook

# Rename worksheets (Excel truncates sheet names to 31 chars; these are
# shortened by one character each in the target revision).
$wb.Worksheets.Item("Include from Attribution Codes").Name = "Include from Attribution Code"
$wb.Worksheets.Item("Include from Attribution Source").Name = "Include from Attribution Sour"
$wb.Worksheets.Item("Include from Attribution Produc").Name = "Include from Attribution Prod"

# Update the Date value on the Metadata sheet.
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B8").Value = "2021-10-01T15:07:10+00:00"
